$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows of data appended to the bottom of the table (rows 56-58)
$rows = @(
    @{ Row=56; A="242743452"; B="Snape";   C="Severus"; D="";         E=20020414; F="Xenobotany Society";          G=72.5 },
    @{ Row=57; A="224303042"; B="White";   C="Walter";  D="Hartwell"; E=20060510; F="ABW";                         G=71.5 },
    @{ Row=58; A="212954131"; B="Baggins"; C="Frodo";   D="M";        E=20040714; F="The Temporal Anomaly Watch";  G=33   }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row

    # Copy formatting (number formats / styles) from the row directly above,
    # restricted to the used columns A:G, so the new rows match existing ones.
    $ws.Range("A" + ($rowIndex - 1) + ":G" + ($rowIndex - 1)).Copy()
    $ws.Range("A" + $rowIndex + ":G" + $rowIndex).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($rowIndex, 1).Value = $r.A
    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D
    $ws.Cells.Item($rowIndex, 5).Value = $r.E
    $ws.Cells.Item($rowIndex, 6).Value = $r.F
    $ws.Cells.Item($rowIndex, 7).Value = $r.G
}

$excel.CutCopyMode = 0

# Extend the duplicate-values conditional formatting on column A to cover the new rows
$cf = $ws.Range("A2:A55").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("A2:A58"))

# Update view state: scroll so row 18 is at the top and select J50 (near the newly added rows)
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J50").Select()
